# Commit: "Small fix to first slide"
#
# The title textbox on slide 1 reads:
#     IETF Hackathon:<TLS 1.3>
# (with a soft line-break between "IETF Hackathon:" and "<TLS 1.3>").
# Remove the stray angle brackets so it reads "TLS 1.3" instead of
# "<TLS 1.3>", without disturbing the line break, the other run, or any
# run-level formatting.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$oldText = "<TLS 1.3>"
$newText = "TLS 1.3"

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        $fullText = $tr.Text
        $idx = $fullText.IndexOf($oldText)
        if ($idx -ge 0) {
            $sub = $tr.Characters($idx + 1, $oldText.Length)
            $sub.Text = $newText
        }
    }
}
